$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price (D) column cells so numeric-looking strings
# (e.g. "21.40", "0.990") are preserved exactly as text, matching the
# source data which stores these values as text, not numbers.
foreach ($addr in @("D2","D3","D5","D9","D10","D12","D13","D15","D16","D17","D18","D25","D26","D32","D33","D36","D39","D40","D41","D42","D44","D46","D47","D49","D51")) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '26.818.88'
$ws.Range('D3').Value = '1.541.53'
$ws.Range('E3').Value = '  -1.63%  '
$ws.Range('E4').Value = '  +0.18%  '
$ws.Range('D5').Value = '205.73'
$ws.Range('E5').Value = '  -0.37%  '
$ws.Range('E6').Value = '  -0.60%  '
$ws.Range('E7').Value = '  +0.20%  '
$ws.Range('E8').Value = '  -0.47%  '
$ws.Range('D9').Value = '21.40'
$ws.Range('E9').Value = '  -2.77%  '
$ws.Range('D10').Value = '0.0581'
$ws.Range('E10').Value = '  -0.52%  '
$ws.Range('E11').Value = '  -0.97%  '
$ws.Range('D12').Value = '1.761.34'
$ws.Range('E12').Value = '  -1.67%  '
$ws.Range('D13').Value = '1.538.60'
$ws.Range('E13').Value = '  -1.79%  '
$ws.Range('E14').Value = '  -1.48%  '
$ws.Range('D15').Value = '0.510'
$ws.Range('E15').Value = '  -0.89%  '
$ws.Range('D16').Value = '26.813.46'
$ws.Range('E16').Value = '  +0.04%  '
$ws.Range('D17').Value = '61.24'
$ws.Range('E17').Value = '  -0.26%  '
$ws.Range('D18').Value = '214.62'
$ws.Range('E18').Value = '  -0.21%  '
$ws.Range('E19').Value = '  -2.41%  '
$ws.Range('E20').Value = '  +0.76%  '
$ws.Range('E21').Value = '  +0.18%  '
$ws.Range('E22').Value = '  -3.16%  '
$ws.Range('E23').Value = '  -1.37%  '
$ws.Range('D25').Value = '152.62'
$ws.Range('E25').Value = '  -0.44%  '
$ws.Range('D26').Value = '6.59'
$ws.Range('E26').Value = '  -2.19%  '
$ws.Range('E27').Value = '  -0.97%  '
$ws.Range('E28').Value = '  +0.16%  '
$ws.Range('E29').Value = '  -0.79%  '
$ws.Range('E30').Value = '  -1.90%  '
$ws.Range('E31').Value = '  -1.57%  '
$ws.Range('D32').Value = '3.22'
$ws.Range('E32').Value = '  +1.70%  '
$ws.Range('D33').Value = '1.367.49'
$ws.Range('E33').Value = '  -2.11%  '
$ws.Range('E34').Value = '  +0.12%  '
$ws.Range('E35').Value = '  -1.46%  '
$ws.Range('D36').Value = '0.964'
$ws.Range('E36').Value = '  +3.10%  '
$ws.Range('E37').Value = '  -0.05%  '
$ws.Range('E38').Value = '  +1.01%  '
$ws.Range('D39').Value = '0.519'
$ws.Range('E39').Value = '  -1.71%  '
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').Value = '5.77'
$ws.Range('E40').Value = '  +8.42%  '
$ws.Range('B41').Value = 'ARBITRUM'
$ws.Range('C41').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D41').Value = '0.805'
$ws.Range('E41').Value = '  -1.14%  '
$ws.Range('D42').Value = '0.990'
$ws.Range('E42').Value = '  +0.13%  '
$ws.Range('E43').Value = '  +0.97%  '
$ws.Range('D44').Value = '63.10'
$ws.Range('E44').Value = '  -0.37%  '
$ws.Range('E45').Value = '  -3.59%  '
$ws.Range('D46').Value = '1.675.79'
$ws.Range('E46').Value = '  -1.67%  '
$ws.Range('D47').Value = '84.22'
$ws.Range('E47').Value = '  -2.08%  '
$ws.Range('E48').Value = '  +3.53%  '
$ws.Range('D49').Value = '0.0₇0977'
$ws.Range('E49').Value = '  -0.77%  '
$ws.Range('E50').Value = '  +0.13%  '
$ws.Range('D51').Value = '0.0940'
